# Edit the response to be a JSON-style short label:
# change cell C2 from "مقبول لم يستدل " to "مقبول "
# then refresh the (wrap-text) row height and move the active selection to D6,
# matching the saved view state produced by Excel after this edit.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

# Update the status text in C2 (shared string reused automatically).
$ws.Range("C2").Value = "مقبول "

# The cell wraps text; shortening it means the row no longer needs the extra
# height it had for the two-line value, so auto-fit row 2 back down.
$ws.Rows.Item(2).EntireRow.AutoFit()

# Move the selection/active cell to D6, matching the workbook's last saved view.
$ws.Range("D6").Select()
